$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D column) updates: prefix with an apostrophe so Excel stores the value
# as text (many of these look like numbers, e.g. "311.50"), then clear the resulting
# quote-prefix formatting so the cell keeps its original (default, unstyled) look ---
$ws.Range("D2").Value = "'43.034.95"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").Value = "'2.308.95"
$ws.Range("D3").ClearFormats()
$ws.Range("D5").Value = "'311.50"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").Value = "'105.99"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").Value = "'0.627"
$ws.Range("D7").ClearFormats()
$ws.Range("D9").Value = "'0.606"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").Value = "'40.12"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").Value = "'0.0912"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "'8.39"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").Value = "'0.108"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").Value = "'0.992"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").Value = "'15.35"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").Value = "'2.655.25"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").Value = "'2.301.98"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").Value = "'42.856.43"
$ws.Range("D18").ClearFormats()
$ws.Range("D19").Value = "'7.49"
$ws.Range("D19").ClearFormats()
$ws.Range("D21").Value = "'13.53"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").Value = "'73.58"
$ws.Range("D22").ClearFormats()
$ws.Range("D24").Value = "'267.65"
$ws.Range("D24").ClearFormats()
$ws.Range("D26").Value = "'7.86"
$ws.Range("D26").ClearFormats()
$ws.Range("D28").Value = "'10.95"
$ws.Range("D28").ClearFormats()
$ws.Range("D30").Value = "'38.25"
$ws.Range("D30").ClearFormats()
$ws.Range("D31").Value = "'22.31"
$ws.Range("D31").ClearFormats()
$ws.Range("D32").Value = "'165.70"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").Value = "'0.0870"
$ws.Range("D33").ClearFormats()
$ws.Range("D34").Value = "'2.79"
$ws.Range("D34").ClearFormats()
$ws.Range("D38").Value = "'0.0358"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").Value = "'2.80"
$ws.Range("D39").ClearFormats()
$ws.Range("D41").Value = "'105.27"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").Value = "'1.58"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").Value = "'71.19"
$ws.Range("D43").ClearFormats()
$ws.Range("D45").Value = "'1.01"
$ws.Range("D45").ClearFormats()
$ws.Range("D46").Value = "'12.30"
$ws.Range("D46").ClearFormats()
$ws.Range("D49").Value = "'76.15"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").Value = "'8.87"
$ws.Range("D50").ClearFormats()
$ws.Range("D51").Value = "'5.18"
$ws.Range("D51").ClearFormats()

# --- Volume(1h) (E column) updates ---
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("E6").Value = "  +2.66%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("E21").Value = "  -1.63%  "
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("E23").Value = "  -1.13%  "
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  +19.45%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("E29").Value = "  -2.04%  "
$ws.Range("E30").Value = "  +2.35%  "
$ws.Range("E31").Value = "  -1.27%  "
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("E34").Value = "  +8.40%  "
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("E39").Value = "  +3.13%  "
$ws.Range("E40").Value = "  -2.65%  "
$ws.Range("E41").Value = "  +9.85%  "
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("E46").Value = "  -1.50%  "
$ws.Range("E49").Value = "  -4.84%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("E51").Value = "  -2.06%  "

# --- Row re-ordering (B/C/D/E) for coins that swapped rank position ---
# Row 36
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'4.64"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.65%  "

# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.111"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.15%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'111.55"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.50%  "

# Row 48
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "'1.695.26"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.68%  "

